$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.113.13'
$ws.Range('E2').Value = '  +0.30%  '

$ws.Range('D3').Value = '3.811.79'
$ws.Range('E3').Value = '  -0.91%  '

$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '698.26'
$ws.Range('E5').Value = '  -0.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.85'
$ws.Range('E6').Value = '  -0.51%  '

$ws.Range('D7').Value = '3.810.87'
$ws.Range('E7').Value = '  -0.88%  '

$ws.Range('E8').Value = '  +0.08%  '

$ws.Range('E9').Value = '  +0.09%  '

$ws.Range('E10').Value = '  -0.49%  '

$ws.Range('E11').Value = '  +2.46%  '

$ws.Range('E12').Value = '  +1.30%  '

$ws.Range('E13').Value = '  -0.96%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.00'
$ws.Range('E14').Value = '  -1.13%  '

$ws.Range('D15').Value = '4.452.34'
$ws.Range('E15').Value = '  -0.96%  '

$ws.Range('D16').Value = '3.861.98'
$ws.Range('E16').Value = '  -1.91%  '

$ws.Range('D17').Value = '71.064.35'
$ws.Range('E17').Value = '  +0.08%  '

$ws.Range('E18').Value = '  +0.64%  '

$ws.Range('E20').Value = '  -0.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '512.42'
$ws.Range('E21').Value = '  +4.04%  '

$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.716'
$ws.Range('E23').Value = '  +0.08%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.03'
$ws.Range('E24').Value = '  -1.24%  '

$ws.Range('E25').Value = '  -1.30%  '

$ws.Range('D26').Value = '3.963.37'
$ws.Range('E26').Value = '  -0.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.07'
$ws.Range('E27').Value = '  -0.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.41'
$ws.Range('E28').Value = '  -1.14%  '

$ws.Range('E29').Value = '  +0.07%  '

$ws.Range('E30').Value = '  -3.60%  '

$ws.Range('E31').Value = '  -5.04%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.43'
$ws.Range('E32').Value = '  -1.00%  '

$ws.Range('E33').Value = '  -1.12%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.10'
$ws.Range('E34').Value = '  -1.02%  '

$ws.Range('E35').Value = '  -4.01%  '

$ws.Range('E36').Value = '  +0.20%  '

$ws.Range('D37').Value = '3.771.47'
$ws.Range('E37').Value = '  -0.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('E39').Value = '  -2.05%  '

$ws.Range('E40').Value = '  +1.04%  '

$ws.Range('E41').Value = '  -0.39%  '

$ws.Range('E42').Value = '  -0.72%  '

$ws.Range('E43').Value = '  -0.13%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '172.93'
$ws.Range('E45').Value = '  +5.50%  '

$ws.Range('E46').Value = '  -0.05%  '

$ws.Range('E47').Value = '  -0.15%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '49.42'
$ws.Range('E48').Value = '  +1.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '429.67'
$ws.Range('E49').Value = '  +4.22%  '

$ws.Range('E50').Value = '  +0.30%  '

$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.37'
$ws.Range('E51').Value = '  +0.67%  '
